$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.19202866666667
$ws.Range("H2").Value = 72.576086
$ws.Range("I2").Value = 0.0688374849512199
$ws.Range("J2").Value = 0.06883748495121988
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 41.17313128465045
$ws.Range("R2").Value = 370.558181561854
$ws.Range("S2").Value = 0.001449457838060978
$ws.Range("T2").Value = 0.001449457838060978
$ws.Range("G3").Value = 24.19202866666667
$ws.Range("H3").Value = 72.576086
$ws.Range("I3").Value = 0.0688374849512199
$ws.Range("J3").Value = 0.06883748495121988
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 1512.097277603814
$ws.Range("R3").Value = 13608.87549843433
$ws.Range("S3").Value = 0.05323183305590848
$ws.Range("T3").Value = 0.05323183305590848
$ws.Range("G4").Value = 24.19202866666667
$ws.Range("H4").Value = 72.576086
$ws.Range("I4").Value = 0.0688374849512199
$ws.Range("J4").Value = 0.06883748495121988
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 10.496171285578
$ws.Range("R4").Value = 94.46554157020202
$ws.Range("S4").Value = 0.0003695069397158864
$ws.Range("T4").Value = 0.0003695069397158863
$ws.Range("G5").Value = 24.19202866666667
$ws.Range("H5").Value = 72.576086
$ws.Range("I5").Value = 0.0688374849512199
$ws.Range("J5").Value = 0.06883748495121988
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 386.404371603008
$ws.Range("R5").Value = 3477.639344427073
$ws.Range("S5").Value = 0.0136029694027621
$ws.Range("T5").Value = 0.01360296940276209
$ws.Range("G6").Value = 24.19202866666667
$ws.Range("H6").Value = 72.576086
$ws.Range("I6").Value = 0.0688374849512199
$ws.Range("J6").Value = 0.06883748495121988
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 5.218664103925557
$ws.Range("R6").Value = 46.96797693533001
$ws.Range("S6").Value = 0.0001837177147724578
$ws.Range("T6").Value = 0.0001837177147724577
$ws.Range("I7").Value = 0.4042872497468467
$ws.Range("J7").Value = 0.4042872497468467
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 241.8126115783112
$ws.Range("R7").Value = 2176.313504204801
$ws.Range("S7").Value = 0.008512764860438128
$ws.Range("T7").Value = 0.008512764860438126
$ws.Range("I8").Value = 0.4042872497468467
$ws.Range("J8").Value = 0.4042872497468467
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.3126341905199885
$ws.Range("T8").Value = 0.3126341905199885
$ws.Range("I9").Value = 0.4042872497468467
$ws.Range("J9").Value = 0.4042872497468467
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 61.64473070050702
$ws.Range("R9").Value = 554.8025763045631
$ws.Range("S9").Value = 0.002170139489058457
$ws.Range("T9").Value = 0.002170139489058457
$ws.Range("I10").Value = 0.4042872497468467
$ws.Range("J10").Value = 0.4042872497468467
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 2269.379260387552
$ws.Range("R10").Value = 20424.41334348797
$ws.Range("S10").Value = 0.0798911681931769
$ws.Range("T10").Value = 0.0798911681931769
$ws.Range("I11").Value = 0.4042872497468467
$ws.Range("J11").Value = 0.4042872497468467
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 30.64957064343279
$ws.Range("R11").Value = 275.8461357908951
$ws.Range("S11").Value = 0.001078986684184724
$ws.Range("T11").Value = 0.001078986684184724
$ws.Range("G12").Value = 65.630404
$ws.Range("H12").Value = 196.891212
$ws.Range("I12").Value = 0.1867487844836031
$ws.Range("J12").Value = 0.1867487844836031
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 111.6983316029187
$ws.Range("R12").Value = 1005.284984426268
$ws.Range("S12").Value = 0.00393222514753311
$ws.Range("T12").Value = 0.00393222514753311
$ws.Range("G13").Value = 65.630404
$ws.Range("H13").Value = 196.891212
$ws.Range("I13").Value = 0.1867487844836031
$ws.Range("J13").Value = 0.1867487844836031
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 4102.159293204588
$ws.Range("R13").Value = 36919.43363884129
$ws.Range("S13").Value = 0.1444123085854958
$ws.Range("T13").Value = 0.1444123085854958
$ws.Range("G14").Value = 65.630404
$ws.Range("H14").Value = 196.891212
$ws.Range("I14").Value = 0.1867487844836031
$ws.Range("J14").Value = 0.1867487844836031
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 28.474997753076
$ws.Range("R14").Value = 256.2749797776841
$ws.Range("S14").Value = 0.001002433076965211
$ws.Range("T14").Value = 0.001002433076965211
$ws.Range("G15").Value = 65.630404
$ws.Range("H15").Value = 196.891212
$ws.Range("I15").Value = 0.1867487844836031
$ws.Range("J15").Value = 0.1867487844836031
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 1048.274014763136
$ws.Range("R15").Value = 9434.466132868223
$ws.Range("S15").Value = 0.03690341102865129
$ws.Range("T15").Value = 0.03690341102865129
$ws.Range("G16").Value = 65.630404
$ws.Range("H16").Value = 196.891212
$ws.Range("I16").Value = 0.1867487844836031
$ws.Range("J16").Value = 0.1867487844836031
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 14.15768136687333
$ws.Range("R16").Value = 127.41913230186
$ws.Range("S16").Value = 0.0004984066449576725
$ws.Range("T16").Value = 0.0004984066449576725
$ws.Range("G17").Value = 62.22373433333333
$ws.Range("H17").Value = 186.671203
$ws.Range("I17").Value = 0.1770552372766232
$ws.Range("J17").Value = 0.1770552372766232
$ws.Range("M17").Value = 1.701929666666667
$ws.Range("N17").Value = 5.105789
$ws.Range("O17").Value = 0.02105622887134972
$ws.Range("P17").Value = 0.02105622887134972
$ws.Range("Q17").Value = 105.9004194326852
$ws.Range("R17").Value = 953.1037748941669
$ws.Range("S17").Value = 0.003728115598967709
$ws.Range("T17").Value = 0.003728115598967708
$ws.Range("G18").Value = 62.22373433333333
$ws.Range("H18").Value = 186.671203
$ws.Range("I18").Value = 0.1770552372766232
$ws.Range("J18").Value = 0.1770552372766232
$ws.Range("O18").Value = 0.7732971809418951
$ws.Range("P18").Value = 0.7732971809418953
$ws.Range("Q18").Value = 3889.228992912747
$ws.Range("R18").Value = 35003.06093621472
$ws.Range("S18").Value = 0.1369163158570111
$ws.Range("T18").Value = 0.1369163158570111
$ws.Range("G19").Value = 62.22373433333333
$ws.Range("H19").Value = 186.671203
$ws.Range("I19").Value = 0.1770552372766232
$ws.Range("J19").Value = 0.1770552372766232
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4338690000000001
$ws.Range("N19").Value = 1.301607
$ws.Range("O19").Value = 0.005367815805265532
$ws.Range("P19").Value = 0.005367815805265533
$ws.Range("Q19").Value = 26.996949391469
$ws.Range("R19").Value = 242.972544523221
$ws.Range("S19").Value = 0.0009503999010584971
$ws.Range("T19").Value = 0.0009503999010584971
$ws.Range("G20").Value = 62.22373433333333
$ws.Range("H20").Value = 186.671203
$ws.Range("I20").Value = 0.1770552372766232
$ws.Range("J20").Value = 0.1770552372766232
$ws.Range("M20").Value = 15.972384
$ws.Range("N20").Value = 47.917152
$ws.Range("O20").Value = 0.1976099128607259
$ws.Range("P20").Value = 0.1976099128607259
$ws.Range("Q20").Value = 993.8613786859839
$ws.Range("R20").Value = 8944.752408173856
$ws.Range("S20").Value = 0.03498787000976866
$ws.Range("T20").Value = 0.03498787000976866
$ws.Range("G21").Value = 62.22373433333333
$ws.Range("H21").Value = 186.671203
$ws.Range("I21").Value = 0.1770552372766232
$ws.Range("J21").Value = 0.1770552372766232
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2157183333333333
$ws.Range("N21").Value = 0.647155
$ws.Range("O21").Value = 0.002668861520763652
$ws.Range("P21").Value = 0.002668861520763652
$ws.Range("Q21").Value = 13.42280026416278
$ws.Range("R21").Value = 120.805202377465
$ws.Range("S21").Value = 0.000472535909817258
$ws.Range("T21").Value = 0.0004725359098172579
$ws.Range("G22").Value = 57.30924366666667
$ws.Range("H22").Value = 171.927731
$ws.Range("I22").Value = 0.1630712435417071
$ws.Range("J22").Value = 0.1630712435417071
$ws.Range("M22").Value = 1.701929666666667
$ws.Range("N22").Value = 5.105789
$ws.Range("O22").Value = 0.02105622887134972
$ws.Range("P22").Value = 0.02105622887134972
$ws.Range("Q22").Value = 97.53630197052877
$ws.Range("R22").Value = 877.826717734759
$ws.Range("S22").Value = 0.003433665426349795
$ws.Range("T22").Value = 0.003433665426349795
$ws.Range("G23").Value = 57.30924366666667
$ws.Range("H23").Value = 171.927731
$ws.Range("I23").Value = 0.1630712435417071
$ws.Range("J23").Value = 0.1630712435417071
$ws.Range("O23").Value = 0.7732971809418951
$ws.Range("P23").Value = 0.7732971809418953
$ws.Range("Q23").Value = 3582.053928751419
$ws.Range("R23").Value = 32238.48535876277
$ws.Range("S23").Value = 0.1261025329234914
$ws.Range("T23").Value = 0.1261025329234914
$ws.Range("G24").Value = 57.30924366666667
$ws.Range("H24").Value = 171.927731
$ws.Range("I24").Value = 0.1630712435417071
$ws.Range("J24").Value = 0.1630712435417071
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.4338690000000001
$ws.Range("N24").Value = 1.301607
$ws.Range("O24").Value = 0.005367815805265532
$ws.Range("P24").Value = 0.005367815805265533
$ws.Range("Q24").Value = 24.864704240413
$ws.Range("R24").Value = 223.782338163717
$ws.Range("S24").Value = 0.0008753363984674804
$ws.Range("T24").Value = 0.0008753363984674805
$ws.Range("G25").Value = 57.30924366666667
$ws.Range("H25").Value = 171.927731
$ws.Range("I25").Value = 0.1630712435417071
$ws.Range("J25").Value = 0.1630712435417071
$ws.Range("M25").Value = 15.972384
$ws.Range("N25").Value = 47.917152
$ws.Range("O25").Value = 0.1976099128607259
$ws.Range("P25").Value = 0.1976099128607259
$ws.Range("Q25").Value = 915.365246593568
$ws.Range("R25").Value = 8238.287219342112
$ws.Range("S25").Value = 0.03222449422636696
$ws.Range("T25").Value = 0.03222449422636696
$ws.Range("G26").Value = 57.30924366666667
$ws.Range("H26").Value = 171.927731
$ws.Range("I26").Value = 0.1630712435417071
$ws.Range("J26").Value = 0.1630712435417071
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2157183333333333
$ws.Range("N26").Value = 0.647155
$ws.Range("O26").Value = 0.002668861520763652
$ws.Range("P26").Value = 0.002668861520763652
$ws.Range("Q26").Value = 12.36265452836722
$ws.Range("R26").Value = 111.263890755305
$ws.Range("S26").Value = 0.0004352145670315404
$ws.Range("T26").Value = 0.0004352145670315404
